$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "AOS" presentation (previously in column F) is reordered to column E,
# and "Quantum" (previously in column E) moves to column F with an updated
# grade (the fifth presentation's grade, 30, is recorded).
$ws.Range("E1").Value = "AOS"
$ws.Range("F1").Value = "Quantum"
$ws.Range("E2").Value = 27
$ws.Range("F2").Value = 30

# Update the dependent formulas so the F/E term order matches the new
# column order of the presentations.
$ws.Range("C6").Formula = "=SUM(B2*B3,C2*C3,D2*D3,F2*F3,E2*E3)/C7"
$ws.Range("C7").Formula = "=SUM(B3,C3,D3,F3,E3)"
$ws.Range("M6").Formula = "=SUM(C2*C3,D2*D3,F2*F3,E2*E3)/SUM(C3,D3,F3,E3)"

# Force recalculation so stored cached values match the new formulas.
$excel.Calculate()

# Update the active cell/selection left by the author after editing.
$ws.Range("F6").Select()
